# Add a new "CLC CO2 Capture" unit process row (row 57) to the Unit Processes
# sheet, mirroring the existing "Kiln System" row (row 56) for cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 56) down onto the
# new row (row 57) so the new cells get the same styles (text number format
# on most columns, etc.) as the rest of the table.
$ws.Range("A56:I56").Copy()
$ws.Range("A57:I57").PasteSpecial(-4122)

# Fill in the new row's values. Order matters for how new shared strings get
# appended, so set the values in this order: display name, varSheet/calcSheet,
# then the remaining cells.
$ws.Range("C57").Value = "CLC CO2 Capture"
$ws.Range("G57").Value = "CLC Capture"
$ws.Range("I57").Value = "CLC Capture"
$ws.Range("A57").Value = "simple_CLC-capture"
$ws.Range("B57").Value = "CCS"
$ws.Range("D57").Value = "CO2"
$ws.Range("E57").Value = "inflow"
$ws.Range("F57").Value = "data/shared/shared_var.xlsx"
$ws.Range("H57").Value = "data/shared/shared_calcs.xlsx"

# Match the author's final cell selection after adding the row.
$ws.Range("E61").Select()
